$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the day-of-week header cell
$ws.Range("A1").Value = "  Wednesday"

# Update the date cell (B1). A bare Value assignment of a date-shaped
# string gets auto-recognised as a real date serial by the COM layer
# (same as real Excel), so: force Text format, enter the literal
# string, then restore B1's original cell format (border/font, General
# number format) by pasting formats from E1 - a cell that already
# shares B1's original style - so we don't leave B1 on a new style.
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B1").Value = "01/04/2020"
$ws.Range("E1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Swap the names in the rota: Anna <-> Supriti
$ws.Range("A3").Value = "Supriti"
$ws.Range("B8").Value = "Anna"

# Narrow columns E, G, H from 16 to 15 characters. ColumnWidth and the
# stored OOXML column width differ by a constant 5/6 offset in this
# engine, so back that out to land exactly on 15.
$targetWidth = 15 - (5/6)
$ws.Columns.Item(5).ColumnWidth = $targetWidth
$ws.Columns.Item(7).ColumnWidth = $targetWidth
$ws.Columns.Item(8).ColumnWidth = $targetWidth
